$wb = $excel.ActiveWorkbook

# --- Sheet: Metadata ---
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("B5").Value = "Birth and Fetal Death Financial Class"
$ws1.Range("B8").Value = "2024-02-27T09:44:15-05:00"
$ws1.Range("B12").Value = "This valueset contains codes to represent birth and fetal death financial class. This valueset is based on `n[PHVS_BirthAndFetalDeathFinancialClass_NCHS](https://phinvads.cdc.gov/vads/ViewValueSet.action?id=D20CD804-8487-E311-AE2A-0017A477041A). using codes from`nthe [Source of Payment Typology](https://terminology.hl7.org/3.0.0/CodeSystem-SOPT.html) codesystem."
$ws1.Range("B14").ClearContents()

# --- Sheet: Include from Local BFDR Codes -> rename + update concept codes ---
$ws2 = $wb.Worksheets.Item("Include from Local BFDR Codes")
$ws2.Name = "Include from Source of Paymen"

# Concept codes change from symbolic ids (finclass_*) to numeric SOPT codes,
# but they must remain stored as text (not auto-converted to numbers), so
# format the column as text first, then assign the new values.
$ws2.Range("A2:A9").NumberFormat = "@"
$ws2.Range("A2").Value = "33"
$ws2.Range("A3").Value = "2"
$ws2.Range("A4").Value = "99"
$ws2.Range("A5").Value = "38"
$ws2.Range("A6").Value = "5"
$ws2.Range("A7").Value = "81"
$ws2.Range("A8").Value = "311"
$ws2.Range("A9").Value = "9999"
$ws2.Range("B11").Value = "https://nahdo.org/sopt"
